$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: add bottom border separator (A4 empty + border; B4/C4/D4/E4 keep values, add border) ---
$ws.Range("A4:E4").Borders(9).LineStyle = 1

# --- Row 5 ---
$ws.Cells.Item(5,2).Value = 382
$ws.Cells.Item(5,3).Value = ' Ah, yes, I\''ve heard!'
$ws.Cells.Item(5,4).Value = ' Ах, да, я слышал об этом!'
$ws.Cells.Item(5,5).Value = ' Àö, äà, ÿ òìúšàì ïá üóïí!'

# --- Row 6 ---
$ws.Cells.Item(6,2).Value = 385
$ws.Cells.Item(6,3).Value = ' The guild is preparing for an\nexpedition soon?!'
$ws.Cells.Item(6,4).Value = ' Гильдия готовится к\nэкспедиции?!'
$ws.Cells.Item(6,5).Value = ' Ãéìûäéÿ ãïóïâéóòÿ ë\nüëòðåäéøéé?!'
$ws.Rows.Item(6).RowHeight = 21.6

# --- Row 7 ---
$ws.Cells.Item(7,2).Value = 388
$ws.Cells.Item(7,3).Value = ' Please do try to be chosen\nas members! ♪'
$ws.Cells.Item(7,4).Value = ' Пожалуйста, постарайтесь, чтобы\nвас в неё взяли! ♪'
$ws.Cells.Item(7,5).Value = ' Ðïçàìôêòóà, ðïòóàñàêóåòû, œóïáú\nâàò â îåæ âèÿìé! ♪'
$ws.Range("A7:E7").Borders(9).LineStyle = 1
$ws.Rows.Item(7).RowHeight = 21.6

# --- Row 8 ---
$ws.Cells.Item(8,2).Value = 352
$ws.Cells.Item(8,3).Value = ' Ah, yes, I\''ve heard!'
$ws.Cells.Item(8,4).Value = ' Ах, да, я слышал об этом!'
$ws.Cells.Item(8,5).Value = ' Àö, äà, ÿ òìúšàì ïá üóïí!'

# --- Row 9 ---
$ws.Cells.Item(9,2).Value = 355
$ws.Cells.Item(9,3).Value = ' You were chosen for the\nexpedition?!'
$ws.Cells.Item(9,4).Value = ' Вас взяли в экспедицию?!'
$ws.Cells.Item(9,5).Value = ' Âàò âèÿìé â üëòðåäéøéý?!'

# --- Row 10 ---
$ws.Cells.Item(10,2).Value = 358
$ws.Cells.Item(10,3).Value = ' Congratulations! I so hope you\ndo very well! ♪'
$ws.Cells.Item(10,4).Value = ' Поздравляю! Я очень надеюсь,\nчто всё пройдёт хорошо! ♪'
$ws.Cells.Item(10,5).Value = ' Ðïèäñàâìÿý! Ÿ ïœåîû îàäåýòû,\nœóï âòæ ðñïêäæó öïñïšï! ♪'
$ws.Rows.Item(10).RowHeight = 21.6

# --- Sheet view: selection + scroll position ---
$wn = $excel.ActiveWindow
$wn.ScrollRow = 4
$wn.ScrollColumn = 1
$ws.Range("D10").Select()
